# Applies the dated worksheet update: refreshes the header date and
# replaces the division problems in the practice table.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Header date
Replace-Text "2025-08-26 Tuesday" "2025-08-27 Wednesday"

# Row 1
Replace-Text "77÷5=" "18÷2="
Replace-Text "26÷6=" "60÷8="
Replace-Text "32÷2=" "23÷9="
Replace-Text "93÷4=" "83÷5="
Replace-Text "52÷8=" "10÷6="

# Row 2 — the five problems are fully reshuffled (not a simple 1:1 text
# swap per cell), so set each cell directly by position instead of
# relying on Find/Replace matching old text.
$t = $d.Tables.Item(1)
$t.Cell(5, 1).Range.Text = "66÷3="
$t.Cell(5, 2).Range.Text = "64÷3="
$t.Cell(5, 3).Range.Text = "75÷3="
$t.Cell(5, 4).Range.Text = "79÷3="
$t.Cell(5, 5).Range.Text = "64÷5="

# Row 3
Replace-Text "62÷5=" "30÷6="
Replace-Text "20÷7=" "81÷6="
Replace-Text "20÷3=" "16÷6="
Replace-Text "50÷4=" "10÷4="
Replace-Text "59÷7=" "77÷7="

# Row 4
Replace-Text "28÷2=" "55÷6="
Replace-Text "17÷5=" "53÷2="
Replace-Text "87÷3=" "32÷7="
Replace-Text "68÷9=" "45÷7="
Replace-Text "29÷9=" "71÷7="

# Row 5
Replace-Text "59÷6=" "33÷3="
Replace-Text "64÷7=" "12÷4="
Replace-Text "61÷8=" "33÷8="
Replace-Text "67÷7=" "46÷2="
Replace-Text "56÷8=" "77÷8="

Write-Output "done"
